$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the table by two rows for the two new collaborators ---
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 7: Prof. Harish Poptani ---
$ws.Range("A7").Value = "Prof. "
$ws.Range("B7").Value = "Harish"
$ws.Range("C7").Value = "Poptani"
$ws.Range("D7").Value = "H."
$ws.Range("E7").Value = "harish.poptani@liverpool.ac.uk"
$ws.Range("F7").Value = "H.Poptani"
$ws.Range("G7").Value = "Liverpool-Physics"
$ws.Range("H7").Value = "Department of Physics, University of Liverpool, Liverpool L69 7ZE, UK"
$ws.Range("I7").Value = 0

# --- Row 8: Prof. F.T. Parambli ---
$ws.Range("A8").Value = "Prof. "
$ws.Range("E8").Value = "f.t.mada-parambil@liverpool.ac.uk"
$ws.Range("C8").Value = "Parambli"
$ws.Range("D8").Value = "F."
$ws.Range("F8").Value = "F.T.Parambil"
$ws.Range("B8").Value = "FT"
$ws.Range("G8").Value = "Liverpool-Physics"
$ws.Range("H8").Value = "Department of Physics, University of Liverpool, Liverpool L69 7ZE, UK"
$ws.Range("I8").Value = 0

# --- Hyperlink the two new email addresses (gives them the Hyperlink style) ---
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:harish.poptani@liverpool.ac.uk") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:f.t.mada-parambil@liverpool.ac.uk") | Out-Null

# --- Update Timothy Greenshaw's email address (row 2) to include the
#     additional hep.ph.liv.ac.uk alias ---
$ws.Range("E2").Value = "green@liverpool.ac.uk; green@hep.ph.liv.ac.uk "

# --- Update the active selection to match the refreshed view ---
$ws.Activate()
$ws.Range("E2").Select()
